$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 2 for the new match (shifts existing rows down by one)
$ws.Rows(2).Insert()
# Excel's row insert copies the header row's bold/centered formatting down;
# the new data row should be plain, like the other data rows.
$ws.Rows(2).ClearFormats()

# 2. Insert a new column at AG (column 33) for the new "Odd_CS_4-4" market,
#    shifting the rest of the correct-score columns one to the right
$ws.Columns("AG").Insert()

# 3. Write the new/updated values into their final positions.
$ws.Cells.Item(1,33).Value = 'Odd_CS_4-4'
$ws.Cells.Item(1,54).Value = 'Odd_CS_2-3_HT'
$ws.Cells.Item(1,55).Value = 'Odd_CS_3-3_HT'
$ws.Cells.Item(1,56).Value = 'Odd_CS_4-4_HT'
$ws.Cells.Item(2,1).Value = 'O6q50van'
$ws.Cells.Item(2,2).Value = '22/11/2024'
$ws.Cells.Item(2,3).Value = '21:00'
$ws.Cells.Item(2,4).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(2,5).Value = 'Oriente Petrolero'
$ws.Cells.Item(2,6).Value = 'Independiente'
$ws.Cells.Item(2,7).Value = 1.55
$ws.Cells.Item(2,8).Value = 4
$ws.Cells.Item(2,9).Value = 6.25
$ws.Cells.Item(2,10).Value = 2.05
$ws.Cells.Item(2,11).Value = 2.4
$ws.Cells.Item(2,12).Value = 5.5
$ws.Cells.Item(2,13).Value = 1.03
$ws.Cells.Item(2,14).Value = 15
$ws.Cells.Item(2,15).Value = 1.2
$ws.Cells.Item(2,16).Value = 4.33
$ws.Cells.Item(2,17).Value = 1.67
$ws.Cells.Item(2,18).Value = 2.15
$ws.Cells.Item(2,19).Value = 1.3
$ws.Cells.Item(2,20).Value = 3.4
$ws.Cells.Item(2,21).Value = 1.73
$ws.Cells.Item(2,22).Value = 2
$ws.Cells.Item(2,23).Value = 8
$ws.Cells.Item(2,24).Value = 8
$ws.Cells.Item(2,25).Value = 8.5
$ws.Cells.Item(2,26).Value = 11
$ws.Cells.Item(2,27).Value = 12
$ws.Cells.Item(2,28).Value = 23
$ws.Cells.Item(2,29).Value = 13
$ws.Cells.Item(2,30).Value = 7.5
$ws.Cells.Item(2,31).Value = 15
$ws.Cells.Item(2,32).Value = 41
$ws.Cells.Item(2,33).Value = 201
$ws.Cells.Item(2,34).Value = 19
$ws.Cells.Item(2,35).Value = 34
$ws.Cells.Item(2,36).Value = 19
$ws.Cells.Item(2,37).Value = 67
$ws.Cells.Item(2,38).Value = 41
$ws.Cells.Item(2,39).Value = 41
$ws.Cells.Item(2,40).Value = 3.6
$ws.Cells.Item(2,41).Value = 7.5
$ws.Cells.Item(2,42).Value = 17
$ws.Cells.Item(2,43).Value = 21
$ws.Cells.Item(2,44).Value = 41
$ws.Cells.Item(2,45).Value = 101
$ws.Cells.Item(2,46).Value = 3.4
$ws.Cells.Item(2,47).Value = 8
$ws.Cells.Item(2,48).Value = 51
$ws.Cells.Item(2,49).Value = 7.5
$ws.Cells.Item(2,50).Value = 29
$ws.Cells.Item(2,51).Value = 34
$ws.Cells.Item(2,52).Value = 101
$ws.Cells.Item(2,53).Value = 101
$ws.Cells.Item(2,54).Value = 201
$ws.Cells.Item(2,55).Value = 51
$ws.Cells.Item(2,56).Value = 51
$ws.Cells.Item(3,7).Value = 1.57
$ws.Cells.Item(3,8).Value = 3.7
$ws.Cells.Item(3,9).Value = 6.25
$ws.Cells.Item(3,10).Value = 2.2
$ws.Cells.Item(3,14).Value = 8.5
$ws.Cells.Item(3,17).Value = 2.08
$ws.Cells.Item(3,18).Value = 1.73
$ws.Cells.Item(3,23).Value = 6
$ws.Cells.Item(3,24).Value = 7
$ws.Cells.Item(3,25).Value = 8.5
$ws.Cells.Item(3,29).Value = 8.5
$ws.Cells.Item(3,42).Value = 21
$ws.Cells.Item(3,43).Value = 26
$ws.Cells.Item(3,54).Value = 351
$ws.Cells.Item(4,7).Value = 1.5
$ws.Cells.Item(4,8).Value = 3.8
$ws.Cells.Item(4,9).Value = 7.5
$ws.Cells.Item(4,12).Value = 7
$ws.Cells.Item(4,17).Value = 2.05
$ws.Cells.Item(4,18).Value = 1.75
$ws.Cells.Item(4,23).Value = 5.5
$ws.Cells.Item(4,29).Value = 8
$ws.Cells.Item(4,32).Value = 81
$ws.Cells.Item(4,33).Value = 900
$ws.Cells.Item(4,36).Value = 23
$ws.Cells.Item(4,39).Value = 67
$ws.Cells.Item(4,40).Value = 3.25
$ws.Cells.Item(4,41).Value = 7.5
$ws.Cells.Item(4,43).Value = 23
$ws.Cells.Item(4,47).Value = 10
$ws.Cells.Item(4,49).Value = 8
$ws.Cells.Item(4,50).Value = 41
$ws.Cells.Item(4,54).Value = 351
$ws.Cells.Item(4,56).Value = 81

# 4. The insert above pushed what used to be the last column (old "Odd_CS_4-4" /
#    BD) out to a new trailing column (BE). We've already copied its values
#    into the correct spot (AG) and re-arranged the two other shuffled
#    trailing columns (Odd_CS_2-3_HT / Odd_CS_3-3_HT / Odd_CS_4-4_HT), so the
#    left-over trailing column is now redundant and must be removed.
$ws.Columns("BE").Delete()
